# Applies the "Side Quest 4/4 Full pack" edits described by the diff.
$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Hunk 1: merge the two runs of
#   "A szörny horda parancsnoka" + ": Hogy merészelsz így szólni hozzám te söpredék."
# into a single run with the same rPr. A Find/Replace across the run
# boundary (identical formatting on both sides) makes Word coalesce the
# two runs into one.
# ---------------------------------------------------------------------
$target1 = "A szörny horda parancsnoka: Hogy merészelsz így szólni hozzám te söpredék."
$d.Content.Find.Execute($target1, $false, $false, $false, $false, $false, $true, 1, $false, $target1, 2) | Out-Null

# ---------------------------------------------------------------------
# Hunk 2: merge
#   "A szörny horda parancsnoka:" + " Arcátlanságodért feldarabollak és megetetlek a vérfarkasaimmal. A nagy "
# into a single run (trailing space preserved -> xml:space="preserve").
# ---------------------------------------------------------------------
$target2 = "A szörny horda parancsnoka: Arcátlanságodért feldarabollak és megetetlek a vérfarkasaimmal. A nagy "
$d.Content.Find.Execute($target2, $false, $false, $false, $false, $false, $true, 1, $false, $target2, 2) | Out-Null

# ---------------------------------------------------------------------
# Hunk 3: wrap the (second) "Kreeber" run - the one in the
# "Kreeber: GYERE TE KORCS!!" line - with a _GoBack bookmark, placed
# before the spellStart proof-error mark and after the run (before
# spellEnd). Rebuild that whole paragraph via InsertXML so the markup
# order matches exactly.
# ---------------------------------------------------------------------
$paragraphs = $d.Paragraphs
for ($i = 1; $i -le $paragraphs.Count; $i++) {
    $candidate = $paragraphs.Item($i)
    if ($candidate.Range.Text -like "*GYERE TE KORCS*") {
        $kreeberPara = $candidate
        break
    }
}
$kreeberXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00646228" w:rsidRDefault="00646228" w:rsidP="004C67B7">
  <w:pPr>
    <w:rPr>
      <w:color w:val="000000" w:themeColor="text1"/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
  </w:pPr>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:color w:val="000000" w:themeColor="text1"/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t>Kreeber</w:t>
  </w:r>
  <w:bookmarkEnd w:id="0"/>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:color w:val="000000" w:themeColor="text1"/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t>: GYERE TE KORCS</w:t>
  </w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      <w:color w:val="000000" w:themeColor="text1"/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t>!!</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
</w:p>
'@
$kreeberPara.Range.InsertXML($kreeberXml) | Out-Null

# ---------------------------------------------------------------------
# Hunk 4: the 5-item quest checklist becomes a 4-item checklist with
# new wording (and the old trailing _GoBack bookmark disappears along
# with the last item).
# ---------------------------------------------------------------------

# 4a. "Menj el a Major épületébe" -> "Beszélj a Majorral"
$d.Content.Find.Execute("Menj el a Major épületébe", $false, $false, $false, $false, $false, $true, 1, $false, "Beszélj a Majorral", 2) | Out-Null

# 4b. The (now second) "Beszélj a Majorral" item becomes three runs:
#     "Menj " + (gramStart) "át  a" (gramEnd) + " portálon"
$paragraphs = $d.Paragraphs
$majorralHits = @()
for ($i = 1; $i -le $paragraphs.Count; $i++) {
    $candidate = $paragraphs.Item($i)
    if ($candidate.Range.Text -eq "Beszélj a Majorral`r") {
        $majorralHits += $i
    }
}
$secondMajorral = $paragraphs.Item($majorralHits[1])
$portalXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00987EBE" w:rsidRDefault="00987EBE" w:rsidP="00987EBE">
  <w:pPr>
    <w:pStyle w:val="Listaszerbekezds"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
    <w:rPr>
      <w:color w:val="000000" w:themeColor="text1"/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:color w:val="000000" w:themeColor="text1"/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t xml:space="preserve">Menj </w:t>
  </w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      <w:color w:val="000000" w:themeColor="text1"/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t>át  a</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:rPr>
      <w:color w:val="000000" w:themeColor="text1"/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t xml:space="preserve"> portálon</w:t>
  </w:r>
</w:p>
'@
$secondMajorral.Range.InsertXML($portalXml) | Out-Null

# 4c. "Keresd meg a szörnyet" -> "Öld meg a szörnyet"
$d.Content.Find.Execute("Keresd meg a szörnyet", $false, $false, $false, $false, $false, $true, 1, $false, "Öld meg a szörnyet", 2) | Out-Null

# 4d. The (now second) "Öld meg a szörnyet" item ->
#     "Menj vissza a portál segítségével és beszélj a Majorral"
$paragraphs = $d.Paragraphs
$szornyetHits = @()
for ($i = 1; $i -le $paragraphs.Count; $i++) {
    $candidate = $paragraphs.Item($i)
    if ($candidate.Range.Text -eq "Öld meg a szörnyet`r") {
        $szornyetHits += $i
    }
}
$secondSzornyet = $paragraphs.Item($szornyetHits[1])
$secondSzornyet.Range.Text = "Menj vissza a portál segítségével és beszélj a Majorral"

# 4e. Delete the trailing "Menj vissza a Major épületébe" item entirely
# (together with its _GoBack bookmark and its own paragraph mark).
$paragraphs = $d.Paragraphs
for ($i = 1; $i -le $paragraphs.Count; $i++) {
    $candidate = $paragraphs.Item($i)
    if ($candidate.Range.Text -like "*Menj vissza a Major épületébe*") {
        $lastItem = $candidate
        break
    }
}
$fullRange = $d.Range($lastItem.Range.Start, $lastItem.Range.End)
$fullRange.Delete() | Out-Null
